# Automatic update of files.
# Applies a row-data permutation: row 4 swaps its species-record data with
# row 8, and row 6 swaps its species-record data with row 7 (columns
# A, B, D, E, F, G, H, Q, R). All other columns are identical between the
# swapped row-pairs, so only these columns need to be written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($ws, $r) {
    $data = @{}
    $data["A"] = $ws.Cells.Item($r, 1).Value2
    $data["B"] = $ws.Cells.Item($r, 2).Value2
    $data["D"] = $ws.Cells.Item($r, 4).Value2
    $data["E"] = $ws.Cells.Item($r, 5).Value2
    $data["F"] = $ws.Cells.Item($r, 6).Value2
    $data["G"] = $ws.Cells.Item($r, 7).Value2
    $data["H"] = $ws.Cells.Item($r, 8).Value2
    $data["Q"] = $ws.Cells.Item($r, 17).Value2
    $data["R"] = $ws.Cells.Item($r, 18).Value2
    return $data
}

function Set-RowData($ws, $r, $data) {
    $ws.Cells.Item($r, 1).Value2 = $data["A"]
    $ws.Cells.Item($r, 2).Value2 = $data["B"]
    $ws.Cells.Item($r, 4).Value2 = $data["D"]
    $ws.Cells.Item($r, 5).Value2 = $data["E"]
    $ws.Cells.Item($r, 6).Value2 = $data["F"]
    $ws.Cells.Item($r, 7).Value2 = $data["G"]
    $ws.Cells.Item($r, 8).Value2 = $data["H"]
    $ws.Cells.Item($r, 17).Value2 = $data["Q"]
    $ws.Cells.Item($r, 18).Value2 = $data["R"]
}

# Capture the original values of the four affected rows before overwriting.
$row4 = Get-RowData $ws 4
$row6 = Get-RowData $ws 6
$row7 = Get-RowData $ws 7
$row8 = Get-RowData $ws 8

# Swap row 4 <-> row 8
Set-RowData $ws 4 $row8
Set-RowData $ws 8 $row4

# Swap row 6 <-> row 7
Set-RowData $ws 6 $row7
Set-RowData $ws 7 $row6

Write-Output "Row swap complete"
